# Update the CDA Logical model ValueSet workbook for ST.r2b:
#  - bump Version / Date metadata values
#  - insert a new "Jurisdiction" property row into the Metadata table
#  - rename the "Include from RoleClass" sheet to "Include #0"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from RoleClass")

# --- Make room for a new "Jurisdiction" row right after "Contact" (row 10) ---
# Capture the existing Property/Value pairs for rows 11-14 before they move.
$a11 = $ws1.Range("A11").Value()
$b11 = $ws1.Range("B11").Value()
$a12 = $ws1.Range("A12").Value()
$b12 = $ws1.Range("B12").Value()
$a13 = $ws1.Range("A13").Value()
$b13 = $ws1.Range("B13").Value()
$a14 = $ws1.Range("A14").Value()
$b14 = $ws1.Range("B14").Value()

# Extend the table formatting down into the new row 15 (copy row 14's look).
$ws1.Range("A14:B14").Copy()
$ws1.Range("A15:B15").PasteSpecial(-4122)

# Shift the captured rows down by one (bottom-up so nothing is clobbered).
$ws1.Range("A15").Value = $a14
$ws1.Range("B15").Value = $b14
$ws1.Range("A14").Value = $a13
$ws1.Range("B14").Value = $b13
$ws1.Range("A13").Value = $a12
$ws1.Range("B13").Value = $b12
$ws1.Range("A12").Value = $a11
$ws1.Range("B12").Value = $b11

# New row 11: Jurisdiction property with no value yet recorded.
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# --- Refresh the Version and Date metadata values ---
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- Rename the include sheet ---
$ws2.Name = "Include #0"
